$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its textual formatting (values like "1.002" or
# "23.322.48" must not be auto-converted to numbers/dates by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: update D2, E2
$ws.Range("D2").Value = '23.322.48'
$ws.Range("E2").Value = '  -0.51%  '

# Row 3: update D3, E3
$ws.Range("D3").Value = '1.624.81'
$ws.Range("E3").Value = '  -0.24%  '

# Row 4: update E4
$ws.Range("E4").Value = '  +0.62%  '

# Row 5: update D5, E5
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.47%  '

# Row 6: update D6, E6
$ws.Range("D6").Value = '303.59'
$ws.Range("E6").Value = '  -1.30%  '

# Row 7: update D7
$ws.Range("D7").Value = '0.3776'

# Row 8: update D8, E8
$ws.Range("D8").Value = '51.83'
$ws.Range("E8").Value = '  -2.26%  '

# Row 9: update D9, E9
$ws.Range("D9").Value = '0.3609'
$ws.Range("E9").Value = '  -1.53%  '

# Row 10: update D10, E10
$ws.Range("D10").Value = '1.230'
$ws.Range("E10").Value = '  -4.14%  '

# Row 11: update D11, E11
$ws.Range("D11").Value = '0.08066'
$ws.Range("E11").Value = '  -1.64%  '

# Row 12: update E12
$ws.Range("E12").Value = '  +0.62%  '

# Row 13: update D13, E13
$ws.Range("D13").Value = '22.56'
$ws.Range("E13").Value = '  -3.03%  '

# Row 14: update D14, E14
$ws.Range("D14").Value = '6.541'
$ws.Range("E14").Value = '  -2.03%  '

# Row 15: update D15, E15
$ws.Range("D15").Value = '0.00001241'
$ws.Range("E15").Value = '  -2.05%  '

# Row 16: update D16, E16
$ws.Range("D16").Value = '7.212'
$ws.Range("E16").Value = '  -3.45%  '

# Row 17: update D17, E17
$ws.Range("D17").Value = '1.625.26'
$ws.Range("E17").Value = '  -0.20%  '

# Row 18: update D18, E18
$ws.Range("D18").Value = '93.36'
$ws.Range("E18").Value = '  -1.57%  '

# Row 19: update D19, E19
$ws.Range("D19").Value = '0.06915'
$ws.Range("E19").Value = '  -0.58%  '

# Row 20: update D20, E20
$ws.Range("D20").Value = '17.90'
$ws.Range("E20").Value = '  -2.91%  '

# Row 21: update D21, E21
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.36%  '

# Row 22: update D22, E22
$ws.Range("D22").Value = '6.436'
$ws.Range("E22").Value = '  -2.37%  '

# Row 23: update D23, E23
$ws.Range("D23").Value = '23.328.97'
$ws.Range("E23").Value = '  -0.51%  '

# Row 24: update D24, E24
$ws.Range("D24").Value = '12.72'
$ws.Range("E24").Value = '  -2.22%  '

# Row 25: update D25, E25
$ws.Range("D25").Value = '3.200'
$ws.Range("E25").Value = '  +2.16%  '

# Row 26: update D26, E26
$ws.Range("D26").Value = '2.446'
$ws.Range("E26").Value = '  +0.74%  '

# Row 27: update D27, E27
$ws.Range("D27").Value = '21.08'
$ws.Range("E27").Value = '  -1.59%  '

# Row 28: update D28, E28
$ws.Range("D28").Value = '148.74'
$ws.Range("E28").Value = '  -1.28%  '

# Row 29: update D29, E29
$ws.Range("D29").Value = '5.284'
$ws.Range("E29").Value = '  +0.03%  '

# Row 30: update D30, E30
$ws.Range("D30").Value = '134.54'
$ws.Range("E30").Value = '  -1.38%  '

# Row 31: update D31, E31
$ws.Range("D31").Value = '2.296'
$ws.Range("E31").Value = '  -5.04%  '

# Row 32: update D32, E32
$ws.Range("D32").Value = '1.807.06'
$ws.Range("E32").Value = '  +0.22%  '

# Row 33: update D33, E33
$ws.Range("D33").Value = '6.728'
$ws.Range("E33").Value = '  -3.22%  '

# Row 34: update D34, E34
$ws.Range("D34").Value = '10.87'
$ws.Range("E34").Value = '  +3.58%  '

# Row 35: update D35, E35
$ws.Range("D35").Value = '0.9444'
$ws.Range("E35").Value = '  -3.24%  '

# Row 36: update D36, E36
$ws.Range("D36").Value = '0.02807'
$ws.Range("E36").Value = '  -0.03%  '

# Row 37: update D37, E37
$ws.Range("D37").Value = '0.2522'
$ws.Range("E37").Value = '  -0.63%  '

# Row 38: update D38, E38
$ws.Range("D38").Value = '0.08809'
$ws.Range("E38").Value = '  -0.59%  '

# Row 39: update D39, E39
$ws.Range("D39").Value = '6.093'
$ws.Range("E39").Value = '  -2.57%  '

# Row 40: update D40, E40
$ws.Range("D40").Value = '0.07092'
$ws.Range("E40").Value = '  -5.32%  '

# Row 41: update D41, E41
$ws.Range("D41").Value = '1.360'
$ws.Range("E41").Value = '  -3.70%  '

# Row 42: update D42, E42
$ws.Range("D42").Value = '0.7020'
$ws.Range("E42").Value = '  -2.26%  '

# Row 43: update D43, E43
$ws.Range("D43").Value = '16.09'
$ws.Range("E43").Value = '  +0.23%  '

# Row 44: update D44, E44
$ws.Range("D44").Value = '12.24'
$ws.Range("E44").Value = '  -4.46%  '

# Row 45: update D45, E45
$ws.Range("D45").Value = '0.6430'
$ws.Range("E45").Value = '  -3.00%  '

# Row 46: update B46, C46, D46, E46
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.44%  '

# Row 47: update B47, C47, D47, E47
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.309'
$ws.Range("E47").Value = '  -2.37%  '

# Row 48: update D48, E48
$ws.Range("D48").Value = '3.980'
$ws.Range("E48").Value = '  -1.49%  '

# Row 49: update D49, E49
$ws.Range("D49").Value = '0.07966'
$ws.Range("E49").Value = '  -0.67%  '

# Row 50: update D50, E50
$ws.Range("D50").Value = '1.199'
$ws.Range("E50").Value = '  -1.39%  '

# Row 51: update D51, E51
$ws.Range("D51").Value = '125.58'
$ws.Range("E51").Value = '  -4.83%  '
